# Problem 3: Gun Violence
#
# Prepare Sheet1 for the new data block: widen the label/value columns
# (B, C, D) that sit next to the already-wide columns (E, F), and leave
# the sheet with column A selected (as happens when a user clicks the
# column-A header before inserting/formatting a new section).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target stored widths (character units) are 30.33203125 / 22.5 / 25.83203125.
# ColumnWidth is quantized to the nearest 1/6 of a character before a fixed
# 5/6-character pad is applied on save, so these inputs land on the closest
# representable widths to the targets (exact for column C).
$ws.Columns.Item(2).ColumnWidth = 29.5                 # column B -> stored width ~30.33
$ws.Columns.Item(3).ColumnWidth = 21.666666666666668   # column C -> stored width 22.5
$ws.Columns.Item(4).ColumnWidth = 25                   # column D -> stored width ~25.83

# Select the entire column A, mirroring a header click (sqref A1:A1048576).
$ws.Columns.Item(1).Select()
